$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '60.543.35'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +3.55%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.645.42'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +0.83%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '570.64'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +6.98%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '147.36'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +3.13%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.995'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -0.38%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.606'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +6.79%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '6.84'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -1.60%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.105'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +4.36%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.144'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +6.63%  '
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +2.86%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.112.67'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +0.67%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '60.521.99'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +3.60%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '21.83'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +5.37%  '
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +4.25%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.657.29'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +1.43%  '
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +3.80%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '345.68'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +3.32%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '10.45'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  +3.29%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.42'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +3.30%  '
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +1.21%  '
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '66.84'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +0.79%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.444'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +6.99%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.996'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.28%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.38'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +4.04%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0₃0784'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +6.65%  '
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -0.08%  '
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +4.51%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.12'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +4.48%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '155.72'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +3.57%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '19.23'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +2.51%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.11'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +5.51%  '
$ws.Range('B36').NumberFormat = "@"
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').NumberFormat = "@"
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.18'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +7.72%  '
$ws.Range('B37').NumberFormat = "@"
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').NumberFormat = "@"
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.918'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +13.59%  '
$ws.Range('B38').NumberFormat = "@"
$ws.Range('B38').Value = 'SuiNetwork'
$ws.Range('C38').NumberFormat = "@"
$ws.Range('C38').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.914'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +7.03%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '37.68'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +1.37%  '
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +7.74%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '309.51'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +10.30%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.68'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +3.16%  '
$ws.Range('B43').NumberFormat = "@"
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').NumberFormat = "@"
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.610'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +2.24%  '
$ws.Range('B44').NumberFormat = "@"
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').NumberFormat = "@"
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.994'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -0.60%  '
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +4.80%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0552'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +4.34%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '19.55'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +2.98%  '
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '125.48'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +11.41%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.975.51'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +1.50%  '
